# alert message for invalid file type
# Appends new log rows recorded while handling an "invalid file type" alert.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$logRows = @(
    @("2024 March 29 9:09:59 AM", "yao", "ReadyImage Signal Recieved"),
    @("2024 March 29 9:10:04 AM", "yao", "ReadyPost Signal Recieved"),
    @("2024 March 29 9:10:18 AM", "yao", "Acc Signal Recieved"),
    @("2024 March 29 9:11:09 AM", "yao", "Image Signal Recieved"),
    @("2024 March 29 9:11:11 AM", "yao", "Post Signal Recieved"),
    @("2024 March 29 9:12:40 AM", "yao", "ReadyImage Signal Recieved"),
    @("2024 March 29 9:12:46 AM", "yao", "ReadyPost Signal Recieved"),
    @("2024 March 29 9:12:55 AM", "yao", "Auth Signal Recieved")
)

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 0 }
$nextRow = $lastRow + 1

foreach ($entry in $logRows) {
    $ws.Cells.Item($nextRow, 1).Value = $entry[0]
    $ws.Cells.Item($nextRow, 2).Value = $entry[1]
    $ws.Cells.Item($nextRow, 3).Value = $entry[2]
    $nextRow++
}
